$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1. Rename the defined name placette.plot_block -> placette.block_code
#    (keeps referring to the same range on the "placette" sheet)
# -----------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "placette.plot_block") {
        $n.Name = "placette.block_code"
    }
}

# -----------------------------------------------------------------
# 2. Rename the "placette" table column / header cell
#    plot_block -> block_code (this also renames the ListColumn)
# -----------------------------------------------------------------
$wsPlacette = $wb.Worksheets.Item("placette")
$wsPlacette.Range('C1').Value = 'block_code'

# -----------------------------------------------------------------
# 3. Update the "dictionary" sheet: the data dictionary is sorted
#    alphabetically by term name (column A). Renaming plot_block to
#    block_code moves its entry from its old alphabetical slot (it
#    used to sit right after "plant_spacing") up to a new slot right
#    after "bbch_stage" (row 3), which shifts every entry that used
#    to sit between "bbch_stage" and "plant_spacing" down by one row.
# -----------------------------------------------------------------
$wsDict = $wb.Worksheets.Item("dictionary")

# Row 3
$wsDict.Range('A3').Value = 'block_code'
$wsDict.Range('B3').Value = 'Code du boc sur la quelle est située l''unité expérimentale. Un bloc regroupe plusieurs unités expérimentales. Les blocs sont utilisés pour réduire la variabilité expérimentale et améliorer la précision des comparaisons entre les traitements.'
$wsDict.Range('G3').Value = 'vignevin:block_code'

# Row 4
$wsDict.Range('A4').Value = 'commune_insee_id'
$wsDict.Range('B4').Value = 'Code INSEE de la commune (5 caractères alphanumériques)'
$wsDict.Range('G4').Value = 'vignevin:commune_insee_id'

# Row 5
$wsDict.Range('A5').Value = 'commune_name'
$wsDict.Range('B5').Value = 'Nom de la commune sur laquelle se trouve la parcelle'
$wsDict.Range('G5').Value = 'vignevin:commune_name'

# Row 6
$wsDict.Range('A6').Value = 'cultivar_name'
$wsDict.Range('B6').Value = 'Nom de la variété (et clone si connu) produisant les fruits. Format de type "Syrah N Cl300" ou "Grenache B". Utiliser la nomenclature de https://www.plantgrape.fr/fr'
$wsDict.Range('G6').Value = 'vignevin:cultivar_name'

# Row 7
$wsDict.Range('A7').Value = 'design_plan'
$wsDict.Range('B7').Value = 'Type de plan d''expérience'
$wsDict.Range('G7').Value = 'vignevin:design_plan'

# Row 8
$wsDict.Range('A8').Value = 'email'
$wsDict.Range('B8').Value = 'Email de la personne'
$wsDict.Range('G8').Value = 'vignevin:email'

# Row 9
$wsDict.Range('A9').Value = 'expe_desc'
$wsDict.Range('B9').Value = 'Description de l’expérimentation et des objectifs poursuivis'
$wsDict.Range('F9').Value = 'character'
$wsDict.Range('G9').Value = 'vignevin:expe_desc'

# Row 10
$wsDict.Range('A10').Value = 'expe_end_date'
$wsDict.Range('B10').Value = 'Date de fin de l’expérimentation. Elle est exprimée au format AAAA-MM-JJ suivant la norme internationale ISO 8601.'
$wsDict.Range('F10').Value = 'date'
$wsDict.Range('G10').Value = 'vignevin:expe_end_date'

# Row 11
$wsDict.Range('A11').Value = 'expe_name'
$wsDict.Range('B11').Value = 'Nom (ou code) usuel de l’expérimentation'
$wsDict.Range('F11').Value = 'character'
$wsDict.Range('G11').Value = 'vignevin:expe_name'

# Row 12
$wsDict.Range('A12').Value = 'expe_start_date'
$wsDict.Range('B12').Value = 'Date de début de l’expérimentation. Elle est exprimée au format AAAA-MM-JJ suivant la norme internationale ISO 8601.'
$wsDict.Range('E12').Value = ''
$wsDict.Range('F12').Value = 'date'
$wsDict.Range('G12').Value = 'vignevin:expe_start_date'

# Row 13
$wsDict.Range('A13').Value = 'field_latitude'
$wsDict.Range('B13').Value = 'Latitude du centroïde de la parcelle  (degrés décimaux WGS84)'
$wsDict.Range('G13').Value = 'vignevin:field_latitude'

# Row 14
$wsDict.Range('A14').Value = 'field_longitude'
$wsDict.Range('B14').Value = 'Longitude du centroïde de la parcelle (degrés décimaux WGS84)'
$wsDict.Range('E14').Value = 'decimal degrees'
$wsDict.Range('F14').Value = 'numeric'
$wsDict.Range('G14').Value = 'vignevin:field_longitude'

# Row 15
$wsDict.Range('A15').Value = 'field_name'
$wsDict.Range('B15').Value = 'Nom (ou code) de la parcelle sur laquelle l’expérimentation a lieu'
$wsDict.Range('F15').Value = 'character'
$wsDict.Range('G15').Value = 'vignevin:field_name'

# Row 16
$wsDict.Range('A16').Value = 'observation_date'
$wsDict.Range('B16').Value = 'Date de réalisation de l''observation.'
$wsDict.Range('F16').Value = 'date'
$wsDict.Range('G16').Value = 'vignevin:observation_date'

# Row 17
$wsDict.Range('A17').Value = 'organization_name'
$wsDict.Range('B17').Value = 'Nom de l''institution responsable de l''expérimentation'
$wsDict.Range('E17').Value = ''
$wsDict.Range('F17').Value = 'character'
$wsDict.Range('G17').Value = 'vignevin:organization_name'

# Row 18
$wsDict.Range('A18').Value = 'plant_spacing'
$wsDict.Range('B18').Value = 'Ecartement entre les ceps de vigne sur le rang, en m'
$wsDict.Range('E18').Value = 'm'
$wsDict.Range('F18').Value = 'numeric'
$wsDict.Range('G18').Value = 'vignevin:plant_spacing'
